$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 448.13635
$ws.Range("I2").Value = 338.6154
$ws.Range("K2").Value = 338.6154
$ws.Range("M2").Value = -225.6154
$ws.Range("H19").Value = 1570.5358
$ws.Range("I19").Value = 1265.2307
$ws.Range("J19").Value = 1835.1333
$ws.Range("K19").Value = 1265.2307
$ws.Range("L19").Value = 1835.1333
$ws.Range("M19").Value = -1090.2307
$ws.Range("N19").Value = -2185.1333
$ws.Range("H38").Value = 5733.846
$ws.Range("I38").Value = 5505
$ws.Range("J38").Value = 6000.8335
$ws.Range("K38").Value = 16515
$ws.Range("L38").Value = 18002.5005
$ws.Range("M38").Value = -16143
$ws.Range("N38").Value = -18746.5005
$ws.Range("H86").Value = 6812.375
$ws.Range("I86").Value = 9250
$ws.Range("K86").Value = 9250
$ws.Range("M86").Value = -8127
$ws.Range("H89").Value = 6812.375
$ws.Range("I89").Value = 9250
$ws.Range("K89").Value = 46250
$ws.Range("M89").Value = -40634
$ws.Range("H111").Value = 4990.143
$ws.Range("I111").Value = 5127.1333
$ws.Range("K111").Value = 15381.3999
$ws.Range("M111").Value = -12314.3999

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 1565.4
$ws.Range("I28").Value = 1565.4
$ws.Range("K28").Value = 1565.4
$ws.Range("M28").Value = -1373.4
$ws.Range("H32").Value = 21284260
$ws.Range("J32").Value = 19659.9
$ws.Range("L32").Value = 19659.9
$ws.Range("N32").Value = -20233.9
$ws.Range("H45").Value = 2868
$ws.Range("J45").Value = 4009.3333
$ws.Range("L45").Value = 4009.3333
$ws.Range("N45").Value = -4763.3333
$ws.Range("H61").Value = 47623456
$ws.Range("I61").Value = 58826590
$ws.Range("K61").Value = 58826590
$ws.Range("M61").Value = -58826378
$ws.Range("H63").Value = 3042.96
$ws.Range("I63").Value = 1946.1316
$ws.Range("K63").Value = 1946.1316
$ws.Range("M63").Value = -1260.1316
$ws.Range("H66").Value = 3042.96
$ws.Range("I66").Value = 1946.1316
$ws.Range("K66").Value = 9730.657999999999
$ws.Range("M66").Value = -6298.657999999999
$ws.Range("H74").Value = 50002476
$ws.Range("I74").Value = 58825650
$ws.Range("K74").Value = 58825650
$ws.Range("M74").Value = -58824776
$ws.Range("H77").Value = 50002476
$ws.Range("I77").Value = 58825650
$ws.Range("K77").Value = 294128250
$ws.Range("M77").Value = -294123882
$ws.Range("H97").Value = 1423.9656
$ws.Range("I97").Value = 1513.5
$ws.Range("J97").Value = 994.2
$ws.Range("K97").Value = 1513.5
$ws.Range("L97").Value = 994.2
$ws.Range("M97").Value = -1017.5
$ws.Range("N97").Value = -1986.2
$ws.Range("H99").Value = 1565.4
$ws.Range("I99").Value = 1565.4
$ws.Range("K99").Value = 1565.4
$ws.Range("M99").Value = 1429.6
$ws.Range("H122").Value = 3107.6365
$ws.Range("I122").Value = 2030.2941
$ws.Range("J122").Value = 4252.3125
$ws.Range("K122").Value = 6090.8823
$ws.Range("L122").Value = 12756.9375
$ws.Range("M122").Value = -3640.8823
$ws.Range("N122").Value = -17656.9375
$ws.Range("H132").Value = 71431780
$ws.Range("I132").Value = 3363.5833
$ws.Range("K132").Value = 10090.7499
$ws.Range("M132").Value = -7560.749899999999
$ws.Range("H135").Value = 46603
$ws.Range("J135").Value = 46603
$ws.Range("L135").Value = 46603
$ws.Range("N135").Value = -56743
$ws.Range("H136").Value = 47623456
$ws.Range("I136").Value = 58826590
$ws.Range("K136").Value = 176479770
$ws.Range("M136").Value = -176477220

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3288.5557
$ws.Range("I107").Value = 2784.3845
$ws.Range("K107").Value = 2784.3845
$ws.Range("M107").Value = -864.3845000000001
$ws.Range("H134").Value = 3594.7097
$ws.Range("I134").Value = 3331.2334
$ws.Range("K134").Value = 9993.700199999999
$ws.Range("M134").Value = -7458.700199999999

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 299.29166
$ws.Range("I7").Value = 217.72728
$ws.Range("J7").Value = 368.30768
$ws.Range("K7").Value = 217.72728
$ws.Range("L7").Value = 368.30768
$ws.Range("M7").Value = -104.72728
$ws.Range("N7").Value = -594.30768
$ws.Range("H31").Value = 19611932
$ws.Range("I31").Value = 3188.6858
$ws.Range("K31").Value = 3188.6858
$ws.Range("M31").Value = -2893.6858
$ws.Range("H34").Value = 19611932
$ws.Range("I34").Value = 3188.6858
$ws.Range("K34").Value = 3188.6858
$ws.Range("M34").Value = -2986.6858
$ws.Range("H62").Value = 3906.818
$ws.Range("I62").Value = 3337
$ws.Range("K62").Value = 3337
$ws.Range("M62").Value = -2713
$ws.Range("H65").Value = 3906.818
$ws.Range("I65").Value = 3337
$ws.Range("K65").Value = 16685
$ws.Range("M65").Value = -13565
$ws.Range("H86").Value = 4065.8333
$ws.Range("I86").Value = 3879
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3879
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2756
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4065.8333
$ws.Range("I89").Value = 3879
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 19395
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -13779
$ws.Range("N89").Value = -36232
$ws.Range("H132").Value = 3266.7083
$ws.Range("I132").Value = 2343.85
$ws.Range("K132").Value = 7031.549999999999
$ws.Range("M132").Value = -4501.549999999999
$ws.Range("H134").Value = 1341.25
$ws.Range("I134").Value = 1230.6666
$ws.Range("K134").Value = 3691.9998
$ws.Range("M134").Value = -1156.9998

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 811.25
$ws.Range("I86").Value = 597.5
$ws.Range("J86").Value = 1025
$ws.Range("K86").Value = 1792.5
$ws.Range("L86").Value = 3075
$ws.Range("M86").Value = -606.5
$ws.Range("N86").Value = -5447
$ws.Range("H89").Value = 811.25
$ws.Range("I89").Value = 597.5
$ws.Range("J89").Value = 1025
$ws.Range("K89").Value = 5377.5
$ws.Range("L89").Value = 9225
$ws.Range("M89").Value = 550.5
$ws.Range("N89").Value = -21081
$ws.Range("H113").Value = 2586.4285
$ws.Range("J113").Value = 3467.3635
$ws.Range("L113").Value = 10402.0905
$ws.Range("N113").Value = -14742.0905
$ws.Range("H122").Value = 1174.75
$ws.Range("J122").Value = 750
$ws.Range("L122").Value = 6750
$ws.Range("N122").Value = -11650
$ws.Range("H132").Value = 1906809.9
$ws.Range("I132").Value = 1256
$ws.Range("J132").Value = 3335975.2
$ws.Range("K132").Value = 11304
$ws.Range("L132").Value = 30023776.8
$ws.Range("M132").Value = -8774
$ws.Range("N132").Value = -30028836.8

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 932
$ws.Range("I80").Value = 800
$ws.Range("K80").Value = 800
$ws.Range("M80").Value = 198
$ws.Range("H83").Value = 932
$ws.Range("I83").Value = 800
$ws.Range("K83").Value = 4000
$ws.Range("M83").Value = 992
$ws.Range("H113").Value = 3510.3809
$ws.Range("I113").Value = 2736.1428
$ws.Range("K113").Value = 2736.1428
$ws.Range("M113").Value = -566.1428000000001
$ws.Range("H122").Value = 6386
$ws.Range("I122").Value = 2973.1
$ws.Range("K122").Value = 8919.299999999999
$ws.Range("M122").Value = -6469.299999999999
$ws.Range("H132").Value = 3389
$ws.Range("I132").Value = 3234.3333
$ws.Range("K132").Value = 9702.999899999999
$ws.Range("M132").Value = -7172.999899999999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2245.4482
$ws.Range("I22").Value = 1374.3636
$ws.Range("J22").Value = 2777.7778
$ws.Range("K22").Value = 1374.3636
$ws.Range("L22").Value = 2777.7778
$ws.Range("M22").Value = -1079.3636
$ws.Range("N22").Value = -3367.7778
$ws.Range("H27").Value = 2245.4482
$ws.Range("I27").Value = 1374.3636
$ws.Range("J27").Value = 2777.7778
$ws.Range("K27").Value = 1374.3636
$ws.Range("L27").Value = 2777.7778
$ws.Range("M27").Value = -1267.3636
$ws.Range("N27").Value = -2991.7778
$ws.Range("H136").Value = 2146.5518
$ws.Range("I136").Value = 1566.6666
$ws.Range("J136").Value = 6371.4287
$ws.Range("K136").Value = 4699.9998
$ws.Range("L136").Value = 19114.2861
$ws.Range("M136").Value = -2149.9998
$ws.Range("N136").Value = -24214.2861

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6565.85
$ws.Range("I132").Value = 6543.6313
$ws.Range("K132").Value = 19630.8939
$ws.Range("M132").Value = -17100.8939
